$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previous used range to remove any stale cell content
$ws.Cells.Clear()

$ws.Cells.Item(1,1).Value = 'Scientific name'
$ws.Cells.Item(1,2).Value = 'Rank'
$ws.Cells.Item(1,3).Value = 'Plankton group'
$ws.Cells.Item(1,4).Value = 'Comment'

$ws.Cells.Item(2,1).Value = 'Unicell'
$ws.Cells.Item(2,3).Value = ' Other microalgae'
$ws.Cells.Item(2,4).Value = 'No rank, checks scientific name only. '

$ws.Cells.Item(3,1).Value = 'Flagellates'
$ws.Cells.Item(3,3).Value = ' Other microalgae'
$ws.Cells.Item(3,4).Value = 'No rank, checks scientific name only. '

$ws.Cells.Item(5,1).Value = 'Cyanobacteria'
$ws.Cells.Item(5,2).Value = 'Phylum'
$ws.Cells.Item(5,3).Value = ' Cyanobacteria'

$ws.Cells.Item(6,1).Value = 'Bacillariophyta'
$ws.Cells.Item(6,2).Value = 'Phylum'
$ws.Cells.Item(6,3).Value = ' Diatoms'

$ws.Cells.Item(7,1).Value = 'Haptophyta'
$ws.Cells.Item(7,2).Value = 'Phylum'
$ws.Cells.Item(7,3).Value = ' Other microalgae'

$ws.Cells.Item(8,1).Value = 'Chlorophyta'
$ws.Cells.Item(8,2).Value = 'Phylum'
$ws.Cells.Item(8,3).Value = ' Other microalgae'

$ws.Cells.Item(9,1).Value = 'Glaucophyta'
$ws.Cells.Item(9,2).Value = 'Phylum'
$ws.Cells.Item(9,3).Value = ' Other microalgae'

$ws.Cells.Item(10,1).Value = 'Ciliophora'
$ws.Cells.Item(10,2).Value = 'Phylum'
$ws.Cells.Item(10,3).Value = ' Ciliates'

$ws.Cells.Item(11,1).Value = 'Cercozoa'
$ws.Cells.Item(11,2).Value = 'Phylum'
$ws.Cells.Item(11,3).Value = ' Other protozoa'

$ws.Cells.Item(12,1).Value = 'Protozoa, classes incertae sedis'
$ws.Cells.Item(12,2).Value = 'Phylum'
$ws.Cells.Item(12,3).Value = ' Other protozoa'

$ws.Cells.Item(13,1).Value = 'Flagellates phylum incertae sedis'
$ws.Cells.Item(13,2).Value = 'Phylum'
$ws.Cells.Item(13,3).Value = ' Other microalgae'

$ws.Cells.Item(14,1).Value = 'Eukarotic picoplankton phylum incertae sedis'
$ws.Cells.Item(14,2).Value = 'Phylum'
$ws.Cells.Item(14,3).Value = ' Other microalgae'

$ws.Cells.Item(15,1).Value = 'Unicells phylum incertae sedis'
$ws.Cells.Item(15,2).Value = 'Phylum'
$ws.Cells.Item(15,3).Value = ' Other microalgae'

$ws.Cells.Item(18,1).Value = 'Dinophyceae'
$ws.Cells.Item(18,2).Value = 'Class'
$ws.Cells.Item(18,3).Value = ' Dinoflagellates'

$ws.Cells.Item(19,1).Value = 'Bacillariophyta'
$ws.Cells.Item(19,2).Value = 'Class'
$ws.Cells.Item(19,3).Value = ' Diatoms'

$ws.Cells.Item(20,1).Value = 'Cryptophyceae'
$ws.Cells.Item(20,2).Value = 'Class'
$ws.Cells.Item(20,3).Value = ' Other microalgae'

$ws.Cells.Item(21,1).Value = 'Bolidophyceae'
$ws.Cells.Item(21,2).Value = 'Class'
$ws.Cells.Item(21,3).Value = ' Other microalgae'

$ws.Cells.Item(22,1).Value = 'Chrysophyceae'
$ws.Cells.Item(22,2).Value = 'Class'
$ws.Cells.Item(22,3).Value = ' Other microalgae'

$ws.Cells.Item(23,1).Value = 'Dictyochophyceae'
$ws.Cells.Item(23,2).Value = 'Class'
$ws.Cells.Item(23,3).Value = ' Other microalgae'

$ws.Cells.Item(24,1).Value = 'Eustigmatophyceae'
$ws.Cells.Item(24,2).Value = 'Class'
$ws.Cells.Item(24,3).Value = ' Other microalgae'

$ws.Cells.Item(25,1).Value = 'Pelagophyceae'
$ws.Cells.Item(25,2).Value = 'Class'
$ws.Cells.Item(25,3).Value = ' Other microalgae'

$ws.Cells.Item(26,1).Value = 'Raphidophyceae'
$ws.Cells.Item(26,2).Value = 'Class'
$ws.Cells.Item(26,3).Value = ' Other microalgae'

$ws.Cells.Item(27,1).Value = 'Synurophyceae'
$ws.Cells.Item(27,2).Value = 'Class'
$ws.Cells.Item(27,3).Value = ' Other microalgae'

$ws.Cells.Item(28,1).Value = 'Coleochaetophyceae'
$ws.Cells.Item(28,2).Value = 'Class'
$ws.Cells.Item(28,3).Value = ' Other microalgae'

$ws.Cells.Item(29,1).Value = 'Klebsormidiophyceae'
$ws.Cells.Item(29,2).Value = 'Class'
$ws.Cells.Item(29,3).Value = ' Other microalgae'

$ws.Cells.Item(30,1).Value = 'Mesostigmatophyceae'
$ws.Cells.Item(30,2).Value = 'Class'
$ws.Cells.Item(30,3).Value = ' Other microalgae'

$ws.Cells.Item(31,1).Value = 'Zygnematophyceae'
$ws.Cells.Item(31,2).Value = 'Class'
$ws.Cells.Item(31,3).Value = ' Other microalgae'

$ws.Cells.Item(32,1).Value = 'Euglenophyceae'
$ws.Cells.Item(32,2).Value = 'Class'
$ws.Cells.Item(32,3).Value = ' Other microalgae'

$ws.Cells.Item(33,1).Value = 'Cryptophyta, ordines incertae sedis'
$ws.Cells.Item(33,2).Value = 'Class'
$ws.Cells.Item(33,3).Value = ' Other protozoa'

$ws.Cells.Item(34,1).Value = 'Bicosoecophyceae'
$ws.Cells.Item(34,2).Value = 'Class'
$ws.Cells.Item(34,3).Value = ' Other protozoa'

$ws.Cells.Item(35,1).Value = 'Bodonophyceae'
$ws.Cells.Item(35,2).Value = 'Class'
$ws.Cells.Item(35,3).Value = ' Other protozoa'

$ws.Cells.Item(36,1).Value = 'Heterokontophyta, ordines incertae sedis'
$ws.Cells.Item(36,2).Value = 'Class'
$ws.Cells.Item(36,3).Value = ' Other protozoa'

$ws.Cells.Item(37,1).Value = 'Craspedophyceae'
$ws.Cells.Item(37,2).Value = 'Class'
$ws.Cells.Item(37,3).Value = ' Other protozoa'

$ws.Cells.Item(38,1).Value = 'Ellobiopsea'
$ws.Cells.Item(38,2).Value = 'Class'
$ws.Cells.Item(38,3).Value = ' Other protozoa'

$ws.Range("C43").Select()
